# Add the "wb" (westbound) worksheet with weekday hourly traffic data,
# and rename the original sheet to "eb" (eastbound).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "eb"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "wb"

# --- Header row -----------------------------------------------------
# Write "Weekday" (G1) before "Hour" (A1) so the new shared-string
# entries land in the same order as the source workbook (Weekday, Hour).
$ws2.Range("G1").Value = "Weekday"
$ws2.Range("A1").Value = "Hour"
$ws2.Range("B1").Value = "Monday"
$ws2.Range("C1").Value = "Tuesday"
$ws2.Range("D1").Value = "Wednesday"
$ws2.Range("E1").Value = "Thursday"
$ws2.Range("F1").Value = "Friday"

# --- Column A: hour-of-day counter -----------------------------------
$ws2.Range("A2").Value = 0
$ws2.Range("A3").Formula = "=A2+1"
$ws2.Range("A4:A25").Formula = "=A3+1"

# --- Columns B:F: observed traffic counts per weekday ----------------
$data = @(
    @(11, 13, 16, 20, 27),
    @(7, 11, 11, 7, 10),
    @(4, 5, 6, 5, 7),
    @(4, 2, 3, 5, 2),
    @(7, 8, 7, 9, 4),
    @(22, 27, 22, 24, 18),
    @(83, 94, 72, 86, 77),
    @(186, 154, 179, 174, 146),
    @(160, 166, 161, 178, 196),
    @(143, 138, 147, 131, 155),
    @(148, 130, 148, 168, 149),
    @(211, 201, 219, 214, 247),
    @(226, 208, 235, 239, 279),
    @(209, 243, 225, 256, 328),
    @(271, 289, 268, 283, 334),
    @(399, 432, 400, 415, 401),
    @(521, 536, 549, 506, 494),
    @(566, 611, 581, 560, 519),
    @(349, 390, 389, 387, 329),
    @(226, 224, 268, 238, 194),
    @(168, 178, 176, 173, 182),
    @(115, 108, 138, 133, 152),
    @(50, 60, 76, 78, 88),
    @(23, 27, 24, 44, 53)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    $ws2.Range("B$row").Value = $values[0]
    $ws2.Range("C$row").Value = $values[1]
    $ws2.Range("D$row").Value = $values[2]
    $ws2.Range("E$row").Value = $values[3]
    $ws2.Range("F$row").Value = $values[4]
}

# --- Column G: weekday average ---------------------------------------
$ws2.Range("G2").Formula = "=AVERAGE(B2:F2)"
$ws2.Range("G3:G25").Formula = "=AVERAGE(B3:F3)"

# --- Formatting: reuse the same cell styles as the "eb" sheet --------
$ws1.Range("A2:A25").Copy() | Out-Null
$ws2.Range("A1:A25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws1.Range("B2:B25").Copy() | Out-Null
$ws2.Range("B1:G25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Column widths (closest values reachable through the width model) -
$ws2.Columns.Item(1).ColumnWidth = 8.33
$ws2.Range("B1:G1").EntireColumn.ColumnWidth = 9.83

# --- Sheet view / selection -------------------------------------------
$ws2.Range("A2").Select()
$ws1.Activate()
